$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 22: MIN @ GSW (away), 2025-05-12 ---
$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "MIN"
$ws.Range("C22").Value = "GSW"
$ws.Range("D22").Value = "away"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2025-05-12"
$ws.Range("E22").ClearFormats()
$ws.Range("F22").Value = "240:00"
$ws.Range("G22").Value = 39
$ws.Range("H22").Value = 79
$ws.Range("I22").Value = 0.494
$ws.Range("J22").Value = 16
$ws.Range("K22").Value = 34
$ws.Range("L22").Value = 0.471
$ws.Range("M22").Value = 23
$ws.Range("N22").Value = 26
$ws.Range("O22").Value = 0.885
$ws.Range("P22").Value = 6
$ws.Range("Q22").Value = 32
$ws.Range("R22").Value = 38
$ws.Range("S22").Value = 22
$ws.Range("T22").Value = 9
$ws.Range("U22").Value = 5
$ws.Range("V22").Value = 19
$ws.Range("W22").Value = 23
$ws.Range("X22").Value = 117
$ws.Range("Y22").Value = 7
$ws.Range("Z22").Value = 27
$ws.Range("AA22").Value = 31
$ws.Range("AB22").Value = 39
$ws.Range("AC22").Value = 20
$ws.Range("AD22").Value = "W"

# --- Row 23: GSW @ MIN (home), 2025-05-12 ---
$ws.Range("A23").Value = 21
$ws.Range("B23").Value = "GSW"
$ws.Range("C23").Value = "MIN"
$ws.Range("D23").Value = "home"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2025-05-12"
$ws.Range("E23").ClearFormats()
$ws.Range("F23").Value = "240:00"
$ws.Range("G23").Value = 37
$ws.Range("H23").Value = 85
$ws.Range("I23").Value = 0.435
$ws.Range("J23").Value = 8
$ws.Range("K23").Value = 27
$ws.Range("L23").Value = 0.296
$ws.Range("M23").Value = 28
$ws.Range("N23").Value = 33
$ws.Range("O23").Value = 0.848
$ws.Range("P23").Value = 14
$ws.Range("Q23").Value = 24
$ws.Range("R23").Value = 38
$ws.Range("S23").Value = 18
$ws.Range("T23").Value = 10
$ws.Range("U23").Value = 1
$ws.Range("V23").Value = 15
$ws.Range("W23").Value = 25
$ws.Range("X23").Value = 110
$ws.Range("Y23").Value = -7
$ws.Range("Z23").Value = 28
$ws.Range("AA23").Value = 32
$ws.Range("AB23").Value = 17
$ws.Range("AC23").Value = 33
$ws.Range("AD23").Value = "L"

# Replicate the bold/bordered/centered style used on the other index
# cells in column A (style index 1) onto the two new index cells,
# the same way it is applied to A2:A21, without introducing new styles.
$ws.Range("A21").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("A23").PasteSpecial(-4122)
$excel.CutCopyMode = 0
